# The paragraph about visiting the campus originally read:
#   "...the exact information [_GoBack]I can enter the Berkeley, prospective
#   student events or any meetings...during my visit."
# followed by an empty paragraph, then "Please let me know...".
#
# The edit:
#   1. Expands the text to "...the exact information I can enter the
#      Berkeley after arriving San Francisco International Airport,
#      prospective student events...during my visit."
#   2. Relocates the (hidden) "_GoBack" bookmark out of the middle of the
#      sentence and into the empty paragraph that follows.
#
# We must move the bookmark first - otherwise the subsequent Find/Replace,
# whose match text spans the bookmark's old (collapsed) position, would
# delete the bookmark outright instead of just the text around it.

$d = $word.ActiveDocument

# Step 1: relocate the hidden "_GoBack" bookmark into the next (empty)
# paragraph. Bookmarks.Add re-defines an existing bookmark of the same
# name at the new range.
$b = $d.Bookmarks.Item("_GoBack")

# Locate the paragraph that currently contains the bookmark.
$containingParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($b.Start -ge $p.Range.Start -and $b.Start -le $p.Range.End) {
        $containingParagraph = $p
        break
    }
}

$followingParagraph = $containingParagraph.Next()
$d.Bookmarks.Add("_GoBack", $followingParagraph.Range)

# Step 2: expand "the exact information I can enter the Berkeley, " into
# "the exact information I can enter the Berkeley after arriving San
# Francisco International Airport, " in one Find/Replace so the two runs
# merge into a single run, matching how Word folds adjoining same-format
# text together after an edit.
$find = $d.Content.Find
$find.Execute(
    "the exact information I can enter the Berkeley, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the exact information I can enter the Berkeley after arriving San Francisco International Airport, ",
    2
)
